$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "E ta enserá skuchamentu aktivo, fasilitashon" "E ta enserá skucha aktivo, fasilitashon"
Replace-Text "A-E-K-P ta para pa e 4 pasonan klave pa fasilitá" "A-E-K-P ta para pa e 4 stapnan klave pa fasilitá"
Replace-Text "Aseptando kontribushon di mayornan ta enkurashá" "Aseptá kontribushon di mayornan ta enkurashá"
Replace-Text "E ta mustra ku bo ta skuchando i ku nan opinion ta konta." "E ta mustra ku bo ta skucha i ku nan opinion ta konta."
Replace-Text "Mi no tabata sinti mi kómodo na promé instante pasobra" "Mi no tabata sintí mi kómodo na kumisamentu pasobra"
Replace-Text "Mi por komprondé ku bo lo sinti bo inkómodo na promé instante pasobra" "Mi por komprondé ku bo lo sinti bo inkómodo den kuminsamentu pasobra"
Replace-Text "Eksplorando ta yuda nos komprondé e eksperensia" "Eksplorá ta yuda nos komprondé e eksperensia"
Replace-Text "E ta sostené solushonamentu di problema i ta yuda" "E ta sostené kon solushoná problema i ta yuda"
Replace-Text "Ehèmpelnan di posibel preguntanan di EKSPLORASHON for di un diskushon anterior" "Ehèmpelnan di posibel preguntanan ku por usa pa EKSPLORA for di un diskushon anterior"
Replace-Text "Kon siguimentu di bo yu su guia a laga bo sinti? Kon bo ta kere" "Kon sigui bo yu su guia a lagá bo sinti? Kon bo ta kere"
Replace-Text "Resumen: Dor di Pasa Tempu ku nos yunan i laga nan tuma e liderazgo" "Resumen: Dor di Pasa Tempu Abo ku Bo Yu i laga nan tuma e liderazgo"
Replace-Text "ora nos simplemente nota kiko nan ta hasiendo ku nos atenshon kompletu." "ora nos simplemente nota kiko nan ta hasi ku nos atenshon kompletu."
Replace-Text "Esaki ta bai bèk na loke nos ta siñando den e promé parti" "Esaki ta bai bèk na loke nos a siña den e promé parti"
Replace-Text "situashon difísil promé ku purba esaki na kas ku su yu." "situashon difísil promé ku e purba esaki na kas ku su yu."
